# BurnDownChart.xlsx - "Edit BurnDown & LogBook"
#
# 1. Log Book (Sheet1) gains 4 new rows of tracked work (rows 34-37),
#    previously blank placeholder rows.
# 2. The "Day 4" column (N) gets an explicit 0 for every existing task
#    row (6-33) instead of being left blank.
# 3. The Burn down chart's "Ideal burndown" total formula widens its
#    SUM range to include the new rows.
# 4. Column widths for Story/Task (C/D) are widened to fit the new text,
#    and the view scrolls/selects near the newly entered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in "Day 4" (column N) with 0 for every existing task row ---
$ws.Range("N6:N33").Value = 0

# --- New Log Book rows 34-37 ---
# Row 34: Alasan memilih Android / Alasan Memilih android
$ws.Range("D34").Value = "Alasan Memilih android"
$ws.Range("D35").Value = "Tips Memilih Android"
$ws.Range("C35").Value = "Tips Android"
$ws.Range("C34").Value = "Alasan memilih Android"
$ws.Range("C36").Value = "Tambahan dari Beberapa Aplikasi yang kami sarankan"
$ws.Range("D36").Value = "5 Aplikasi Android untuk menjadikan ponsel sebagai scanner tulisan"
$ws.Range("C37").Value = "Tambahan dari Beberapa Aplikasi yang kami sarankan"
$ws.Range("D37").Value = "Aplikasi terbaik untuk keamanan android dan cara melacak perangkat android yang hilang"

$ws.Range("E34:N34").Value = 3
$ws.Range("E35:N35").Value = 1
$ws.Range("E36:N36").Value = 2
$ws.Range("E37:N37").Value = 2

# --- Widen the "Ideal - Remaining efforts" total to cover the new rows ---
$ws.Range("E45").Formula = "=SUM(E6:E37)"

# --- Column widths (Story / Task) widened to fit the longer new text ---
$ws.Columns.Item(3).ColumnWidth = 48.43
$ws.Columns.Item(4).ColumnWidth = 37.166666666666664

# --- Move the view near the newly entered rows ---
$ws.Range("N37").Select()

$wb.Save()
